$d = $word.ActiveDocument

# Insert a comma after the closing parenthesis in
# "(https://obsproject.com/pt-br/) que grava a tela do computador"
# so it reads "(https://obsproject.com/pt-br/), que grava a tela do computador"
$find = $d.Content.Find
$find.Execute(") que grava a tela do computador", $true, $false, $false, $false, $false,
              $true, 1, $false, "), que grava a tela do computador", 2)
